$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The account-statement database was updated: the existing three periods
# (2507/2506/2505) are re-sorted ascending (2505/2506/2507) and a new
# period (2508) is appended as an additional data row, with the "Cant.
# Periodos" count and the total "Valor Mora" updated accordingly.

# Step 1: insert one new blank row at row 20 (a row with no data/special
# formatting nearby) so everything from row 20 down (the blank rows plus
# the signature-block footer) shifts down by one, without disturbing the
# existing data rows 15-18.
$ws.Rows.Item(20).Insert()

# Step 2: row 19 is now free (still blank, untouched by the insert above).
# Copy row 18 - which carries the worker data plus the special "last row"
# bottom-border styling - down into row 19, so that styling/content moves
# to what is now the last data row.
$ws.Range("B18:J18").Copy($ws.Range("B19:J19"))

# Step 3: copy row 17's styling (the plain "middle row" style, no special
# bottom border) onto row 18, since row 18 is no longer the last row.
$ws.Range("B17:J17").Copy($ws.Range("B18:J18"))

# Step 4: write the correct ascending period values into the four data
# rows.
$ws.Range("E16").Value = "2505"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2507"
$ws.Range("E19").Value = "2508"

# Step 5: refresh the summary figures - total overdue amount and period
# count - to reflect the newly added period.
$ws.Range("E11").Value = 227760
$ws.Range("F13").Value = 4
